$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46029
$ws.Range("B2").Value = 95.59
$ws.Range("C2").Value = 90.34
$ws.Range("D2").Value = 86.52
$ws.Range("E2").Value = 82.38
$ws.Range("F2").Value = 81.95999999999999
$ws.Range("G2").Value = 87.83
$ws.Range("H2").Value = 95.91
$ws.Range("I2").Value = 110.49
$ws.Range("J2").Value = 116.51
$ws.Range("K2").Value = 112.61
$ws.Range("L2").Value = 103.84
$ws.Range("M2").Value = 101.46
$ws.Range("N2").Value = 96.15000000000001
$ws.Range("O2").Value = 94.31
$ws.Range("P2").Value = 90.56999999999999
$ws.Range("Q2").Value = 98
$ws.Range("R2").Value = 106.12
$ws.Range("S2").Value = 115.46
$ws.Range("T2").Value = 134.94
$ws.Range("U2").Value = 132.06
$ws.Range("V2").Value = 119.85
$ws.Range("W2").Value = 109.62
$ws.Range("X2").Value = 103.02
$ws.Range("Y2").Value = 90.33
$ws.Range("Z2").Value = 102.33
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 122.14
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 133.5
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 114.74
$ws.Range("AG2").Value = "0h-23h"
